# Apply cryptos-list price/volume refresh per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.575.17"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").Value = "2.291.21"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'322.57"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").Value = "'104.03"
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "'40.04"
$ws.Range("E10").Value = "  +2.87%  "
$ws.Range("D11").Value = "'0.0907"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "'8.58"
$ws.Range("E12").Value = "  +3.59%  "
$ws.Range("D14").Value = "'0.970"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").Value = "'15.27"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "2.292.51"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "42.387.75"
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("D19").Value = "'7.45"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("E21").Value = "  +34.05%  "
$ws.Range("D22").Value = "'73.31"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").Value = "'3.58"
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("D24").Value = "'269.37"
$ws.Range("E24").Value = "  -5.57%  "
$ws.Range("D25").Value = "'2.22"
$ws.Range("E25").Value = "  -1.68%  "
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "'10.87"
$ws.Range("E28").Value = "  +2.23%  "
$ws.Range("D29").Value = "'22.55"
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("D30").Value = "'38.28"
$ws.Range("E30").Value = "  +10.37%  "
$ws.Range("D31").Value = "'165.56"
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("D32").Value = "'6.16"
$ws.Range("E32").Value = "  +5.79%  "
$ws.Range("D33").Value = "'0.0883"
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'2.53"
$ws.Range("E35").Value = "  -13.09%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.114"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D39").Value = "'3.74"
$ws.Range("E39").Value = "  +3.98%  "
$ws.Range("E40").Value = "  -5.62%  "
$ws.Range("E41").Value = "  +5.29%  "
$ws.Range("D42").Value = "'69.88"
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").Value = "'96.13"
$ws.Range("E43").Value = "  -6.79%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").Value = "'12.40"
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("D47").Value = "'81.08"
$ws.Range("E47").Value = "  +5.67%  "
$ws.Range("D48").Value = "'113.03"
$ws.Range("E48").Value = "  -1.71%  "
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").Value = "1.582.81"
$ws.Range("E51").Value = "  +3.16%  "
